# Add a new worksheet that reports the number of effect sizes / studies
# for each quality-score level of the moderator (per outcome).
$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts the new sheet before the active sheet, so we
# create it and then move it to the end of the workbook (after
# "coefficients") to match the target layout.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "nr_studies"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# Re-fetch the worksheet by name: after Move() the old object reference
# can end up pointing at a different sheet, so look it up fresh.
$ws = $wb.Worksheets.Item("nr_studies")

# ---- Header row ----
$ws.Cells.Item(1, 1).Value = "outcome"
$ws.Cells.Item(1, 2).Value = "quality_score_out_of_5"
$ws.Cells.Item(1, 3).Value = "n_effect_sizes"
$ws.Cells.Item(1, 4).Value = "k_studies"

# Match the bold + centered header style already used on the other
# sheets by copying the formatting from an existing header cell.
$styleSource = $wb.Worksheets.Item("descriptives").Range("A1")
$styleSource.Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# ---- Data rows ----
# quality_score_out_of_5 values are stored as text (e.g. "2"), not numbers.
$data = @(
    @("NS", "2", 64, 8),
    @("NS", "5", 59, 8),
    @("NS", "4", 347, 31),
    @("NS", "3", 243, 31),
    @("NS", "1", 9, 3),
    @("NT", "2", 37, 4),
    @("NT", "5", 51, 8),
    @("NT", "4", 177, 21),
    @("NT", "3", 116, 18)
)

$rowIdx = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIdx, 1).Value = $row[0]

    # Force the quality score to be stored as text rather than a number.
    $qCell = $ws.Cells.Item($rowIdx, 2)
    $qCell.NumberFormat = "@"
    $qCell.Value = $row[1]
    $qCell.Style = "Normal"

    $ws.Cells.Item($rowIdx, 3).Value = $row[2]
    $ws.Cells.Item($rowIdx, 4).Value = $row[3]
    $rowIdx++
}
